$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), row 17
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 1605.5366
$ws.Range("J17").Value = 1608.425
$ws.Range("L17").Value = 4825.275
$ws.Range("N17").Value = -5161.275

# Sheet ALC (index 1), row 87
$ws = $wb.Worksheets.Item(1)
$ws.Range("H87").Value = 53000
$ws.Range("J87").Value = 95000
$ws.Range("L87").Value = 95000
$ws.Range("N87").Value = -97496

# Sheet ALC (index 1), row 90
$ws = $wb.Worksheets.Item(1)
$ws.Range("H90").Value = 53000
$ws.Range("J90").Value = 95000
$ws.Range("L90").Value = 285000
$ws.Range("N90").Value = -297480

# Sheet ALC (index 1), row 116
$ws = $wb.Worksheets.Item(1)
$ws.Range("H116").Value = 4709.3887
$ws.Range("I116").Value = 4835.0835
$ws.Range("K116").Value = 4835.0835
$ws.Range("M116").Value = -1393.0835

# Sheet ALC (index 1), row 118
$ws = $wb.Worksheets.Item(1)
$ws.Range("H118").Value = 756
$ws.Range("J118").Value = 900
$ws.Range("L118").Value = 2700
$ws.Range("N118").Value = -6014

# Sheet ALC (index 1), row 125
$ws = $wb.Worksheets.Item(1)
$ws.Range("H125").Value = 1392.75
$ws.Range("I125").Value = 1035
$ws.Range("J125").Value = 1512
$ws.Range("K125").Value = 9315
$ws.Range("L125").Value = 13608
$ws.Range("M125").Value = -6855
$ws.Range("N125").Value = -18528

# Sheet ALC (index 1), row 127
$ws = $wb.Worksheets.Item(1)
$ws.Range("H127").Value = 470.8
$ws.Range("I127").Value = 470.8
$ws.Range("K127").Value = 1412.4
$ws.Range("M127").Value = 3547.6

# Sheet ALC (index 1), row 131
$ws = $wb.Worksheets.Item(1)
$ws.Range("H131").Value = 835681.4399999999
$ws.Range("I131").Value = 1001917.9
$ws.Range("K131").Value = 3005753.7
$ws.Range("M131").Value = -3000713.7

# Sheet ARM (index 2), row 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1688.4445
$ws.Range("J2").Value = 1466
$ws.Range("L2").Value = 1466
$ws.Range("N2").Value = -1692

# Sheet ARM (index 2), row 110
$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 3509.1667
$ws.Range("I110").Value = 3509.1667
$ws.Range("K110").Value = 3509.1667
$ws.Range("M110").Value = -1464.1667

# Sheet ARM (index 2), row 116
$ws = $wb.Worksheets.Item(2)
$ws.Range("H116").Value = 1688.4445
$ws.Range("J116").Value = 1466
$ws.Range("L116").Value = 1466
$ws.Range("N116").Value = -6054

# Sheet ARM (index 2), row 119
$ws = $wb.Worksheets.Item(2)
$ws.Range("H119").Value = 59500
$ws.Range("J119").Value = 59500
$ws.Range("L119").Value = 59500
$ws.Range("N119").Value = -69176

# Sheet ARM (index 2), row 129
$ws = $wb.Worksheets.Item(2)
$ws.Range("H129").Value = 24375
$ws.Range("J129").Value = 24375
$ws.Range("L129").Value = 24375
$ws.Range("N129").Value = -34375

# Sheet BSM (index 3), row 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1688.4445
$ws.Range("J3").Value = 1466
$ws.Range("L3").Value = 1466
$ws.Range("N3").Value = -1694

# Sheet BSM (index 3), row 80
$ws = $wb.Worksheets.Item(3)
$ws.Range("H80").Value = 315.73685
$ws.Range("I80").Value = 210
$ws.Range("K80").Value = 210
$ws.Range("M80").Value = 788

# Sheet BSM (index 3), row 83
$ws = $wb.Worksheets.Item(3)
$ws.Range("H83").Value = 315.73685
$ws.Range("I83").Value = 210
$ws.Range("K83").Value = 1050
$ws.Range("M83").Value = 3942

# Sheet BSM (index 3), row 86
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 3027.348
$ws.Range("J86").Value = 3444.9167
$ws.Range("L86").Value = 3444.9167
$ws.Range("N86").Value = -5690.9167

# Sheet BSM (index 3), row 89
$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value = 3027.348
$ws.Range("J89").Value = 3444.9167
$ws.Range("L89").Value = 17224.5835
$ws.Range("N89").Value = -28456.5835

# Sheet BSM (index 3), row 107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 5495932.5
$ws.Range("I107").Value = 8548148
$ws.Range("J107").Value = 1943.8
$ws.Range("K107").Value = 8548148
$ws.Range("L107").Value = 1943.8
$ws.Range("M107").Value = -8546228
$ws.Range("N107").Value = -5783.8

# Sheet CRP (index 4), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 6171.516
$ws.Range("I31").Value = 4665.3335
$ws.Range("K31").Value = 4665.3335
$ws.Range("M31").Value = -4370.3335

# Sheet CRP (index 4), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 6171.516
$ws.Range("I34").Value = 4665.3335
$ws.Range("K34").Value = 4665.3335
$ws.Range("M34").Value = -4463.3335

# Sheet CRP (index 4), row 99
$ws = $wb.Worksheets.Item(4)
$ws.Range("H99").Value = 3999.75
$ws.Range("I99").Value = 1999.5
$ws.Range("K99").Value = 1999.5
$ws.Range("M99").Value = -501.5

# Sheet CRP (index 4), row 122
$ws = $wb.Worksheets.Item(4)
$ws.Range("H122").Value = 1544
$ws.Range("I122").Value = 1554.3572
$ws.Range("J122").Value = 1399
$ws.Range("K122").Value = 4663.071599999999
$ws.Range("L122").Value = 4197
$ws.Range("M122").Value = -2213.071599999999
$ws.Range("N122").Value = -9097

# Sheet CRP (index 4), row 126
$ws = $wb.Worksheets.Item(4)
$ws.Range("H126").Value = 3999.75
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5

# Sheet CUL (index 5), row 22
$ws = $wb.Worksheets.Item(5)
$ws.Range("H22").Value = 3295.1875
$ws.Range("I22").Value = 957.5
$ws.Range("J22").Value = 4074.4167
$ws.Range("K22").Value = 2872.5
$ws.Range("L22").Value = 12223.2501
$ws.Range("M22").Value = -2703.5
$ws.Range("N22").Value = -12561.2501

# Sheet CUL (index 5), row 27
$ws = $wb.Worksheets.Item(5)
$ws.Range("H27").Value = 3295.1875
$ws.Range("I27").Value = 957.5
$ws.Range("J27").Value = 4074.4167
$ws.Range("K27").Value = 2872.5
$ws.Range("L27").Value = 12223.2501
$ws.Range("M27").Value = -2770.5
$ws.Range("N27").Value = -12427.2501

# Sheet CUL (index 5), row 59
$ws = $wb.Worksheets.Item(5)
$ws.Range("H59").Value = 1024.75
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1024.75
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 3074.25
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -4154.25

# Sheet CUL (index 5), row 62
$ws = $wb.Worksheets.Item(5)
$ws.Range("H62").Value = 6507
$ws.Range("J62").Value = 6507
$ws.Range("L62").Value = 19521
$ws.Range("N62").Value = -20893

# Sheet CUL (index 5), row 65
$ws = $wb.Worksheets.Item(5)
$ws.Range("H65").Value = 6507
$ws.Range("J65").Value = 6507
$ws.Range("L65").Value = 58563
$ws.Range("N65").Value = -65427

# Sheet CUL (index 5), row 68
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 10007662
$ws.Range("J68").Value = 11119536
$ws.Range("L68").Value = 33358608
$ws.Range("N68").Value = -33360230

# Sheet CUL (index 5), row 71
$ws = $wb.Worksheets.Item(5)
$ws.Range("H71").Value = 10007662
$ws.Range("J71").Value = 11119536
$ws.Range("L71").Value = 100075824
$ws.Range("N71").Value = -100083936

# Sheet CUL (index 5), row 93
$ws = $wb.Worksheets.Item(5)
$ws.Range("H93").Value = 8666.666999999999
$ws.Range("J93").Value = 8666.666999999999
$ws.Range("L93").Value = 26000.001
$ws.Range("N93").Value = -29744.001

# Sheet CUL (index 5), row 95
$ws = $wb.Worksheets.Item(5)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Sheet CUL (index 5), row 122
$ws = $wb.Worksheets.Item(5)
$ws.Range("H122").Value = 1499.8572
$ws.Range("J122").Value = 1594.4736
$ws.Range("L122").Value = 14350.2624
$ws.Range("N122").Value = -19250.2624

# Sheet CUL (index 5), row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 2243.1765
$ws.Range("I131").Value = 3184.5
$ws.Range("K131").Value = 9553.5
$ws.Range("M131").Value = -4513.5

# Sheet GSM (index 6), row 68
$ws = $wb.Worksheets.Item(6)
$ws.Range("H68").Value = 86249.75
$ws.Range("J68").Value = 86249.75
$ws.Range("L68").Value = 86249.75

# Sheet GSM (index 6), row 69
$ws = $wb.Worksheets.Item(6)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Sheet GSM (index 6), row 71
$ws = $wb.Worksheets.Item(6)
$ws.Range("H71").Value = 86249.75
$ws.Range("J71").Value = 86249.75
$ws.Range("L71").Value = 258749.25
$ws.Range("N71").Value = -266861.25

# Sheet GSM (index 6), row 72
$ws = $wb.Worksheets.Item(6)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Sheet GSM (index 6), row 80
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 142862050
$ws.Range("I80").Value = 333337340
$ws.Range("J80").Value = 5600
$ws.Range("K80").Value = 333337340
$ws.Range("L80").Value = 5600
$ws.Range("M80").Value = -333336342
$ws.Range("N80").Value = -7596

# Sheet GSM (index 6), row 83
$ws = $wb.Worksheets.Item(6)
$ws.Range("H83").Value = 142862050
$ws.Range("I83").Value = 333337340
$ws.Range("J83").Value = 5600
$ws.Range("K83").Value = 1666686700
$ws.Range("L83").Value = 28000
$ws.Range("M83").Value = -1666681708
$ws.Range("N83").Value = -37984

# Sheet GSM (index 6), row 122
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 2852168
$ws.Range("I122").Value = 5130873.5
$ws.Range("J122").Value = 3786.0833
$ws.Range("K122").Value = 15392620.5
$ws.Range("L122").Value = 11358.2499
$ws.Range("M122").Value = -15390170.5
$ws.Range("N122").Value = -16258.2499

# Sheet LTW (index 7), row 40
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 3951.6428
$ws.Range("I40").Value = 2302.2727
$ws.Range("J40").Value = 9999.333000000001
$ws.Range("K40").Value = 2302.2727
$ws.Range("L40").Value = 9999.333000000001
$ws.Range("M40").Value = -2166.2727
$ws.Range("N40").Value = -10271.333

# Sheet LTW (index 7), row 121
$ws = $wb.Worksheets.Item(7)
$ws.Range("H121").Value = 49475.6
$ws.Range("J121").Value = 49475.6
$ws.Range("L121").Value = 49475.6
$ws.Range("N121").Value = -52969.6

# Sheet LTW (index 7), row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 11918.708
$ws.Range("I122").Value = 10662.454
$ws.Range("J122").Value = 12981.692
$ws.Range("K122").Value = 31987.362
$ws.Range("L122").Value = 38945.076
$ws.Range("M122").Value = -29537.362
$ws.Range("N122").Value = -43845.076

# Sheet WVR (index 8), row 81
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = 4706.5
$ws.Range("I81").Value = 4412.4165
$ws.Range("J81").Value = 5294.6665
$ws.Range("K81").Value = 8824.833000000001
$ws.Range("L81").Value = 10589.333
$ws.Range("M81").Value = -7763.833000000001
$ws.Range("N81").Value = -12711.333

# Sheet WVR (index 8), row 84
$ws = $wb.Worksheets.Item(8)
$ws.Range("H84").Value = 4706.5
$ws.Range("I84").Value = 4412.4165
$ws.Range("J84").Value = 5294.6665
$ws.Range("K84").Value = 44124.165
$ws.Range("L84").Value = 52946.665
$ws.Range("M84").Value = -38820.165
$ws.Range("N84").Value = -63554.665

# Sheet WVR (index 8), row 107
$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 1272.1852
$ws.Range("I107").Value = 999.4761999999999
$ws.Range("J107").Value = 2226.6667
$ws.Range("K107").Value = 2998.4286
$ws.Range("L107").Value = 6680.000100000001
$ws.Range("M107").Value = -1078.4286
$ws.Range("N107").Value = -10520.0001

# Sheet WVR (index 8), row 119
$ws = $wb.Worksheets.Item(8)
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Sheet WVR (index 8), row 121
$ws = $wb.Worksheets.Item(8)
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

# Sheet WVR (index 8), row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 2204.4546
$ws.Range("I132").Value = 2041.02
$ws.Range("K132").Value = 6123.059999999999
$ws.Range("M132").Value = -3593.059999999999
